$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("K2").Value = "32,62 TL - 32,62 TL"

# Row 3
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 4
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 5
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 6
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = "8.300,01 TL - 199,41 TL"
$ws.Range("E6").Value = ""

# Row 7
$ws.Range("K7").Value = "%3,1"

# Row 8
$ws.Range("C8").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 9
$ws.Range("C9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 10
$ws.Range("C10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# Row 11
$ws.Range("C11").Value = ""
$ws.Range("E11").Value = ""

# Row 12
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = "WU: 0,75 USD–12 USD; Diğer: 700 TL–4.000 TL"
$ws.Range("K12").Value = "WU: ,USD–; Diğer: 404,16 TL–3.403,42 TL"

# Row 13
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = "Hesaba: Asgari 300 TL | Azami 3.080 TL"
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 865,75 TL"

# Row 14
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = "3.500 TL - 13.500 TL"
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = "1.952,38 TL - 9.523,81 TL"
